$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "64.071.54"
Set-TextValue "E2" "  +3.44%  "
Set-TextValue "D3" "3.061.74"
Set-TextValue "E3" "  +2.22%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "561.36"
Set-TextValue "E5" "  +3.10%  "
Set-TextValue "D6" "143.40"
Set-TextValue "E6" "  +3.64%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "3.060.87"
Set-TextValue "E8" "  +2.38%  "
Set-TextValue "E10" "  +5.56%  "
Set-TextValue "D11" "6.11"
Set-TextValue "E11" "  -9.29%  "
Set-TextValue "D12" "0.497"
Set-TextValue "E12" "  +12.57%  "
Set-TextValue "E13" "  +5.81%  "
Set-TextValue "D14" "35.53"
Set-TextValue "E14" "  +5.37%  "
Set-TextValue "D15" "3.561.04"
Set-TextValue "E15" "  +2.31%  "
Set-TextValue "D16" "64.094.02"
Set-TextValue "E16" "  +3.47%  "
Set-TextValue "D17" "3.063.75"
Set-TextValue "E17" "  +2.28%  "
Set-TextValue "E18" "  +2.73%  "
Set-TextValue "D19" "6.78"
Set-TextValue "E19" "  +4.15%  "
Set-TextValue "D20" "478.86"
Set-TextValue "E20" "  +3.55%  "
Set-TextValue "D21" "13.96"
Set-TextValue "E21" "  +5.64%  "
Set-TextValue "D22" "0.683"
Set-TextValue "E22" "  +5.26%  "
Set-TextValue "D23" "14.42"
Set-TextValue "E23" "  +15.51%  "
Set-TextValue "E24" "  +5.62%  "
Set-TextValue "D25" "82.50"
Set-TextValue "E25" "  +4.65%  "
Set-TextValue "E26" "  +0.02%  "
Set-TextValue "E27" "  +4.04%  "
Set-TextValue "D28" "8.13"
Set-TextValue "E28" "  +7.65%  "
Set-TextValue "D29" "2.04"
Set-TextValue "E29" "  +2.60%  "
Set-TextValue "D30" "1.00"
Set-TextValue "E30" "  +0.03%  "
Set-TextValue "D31" "26.33"
Set-TextValue "E31" "  +4.28%  "
Set-TextValue "E32" "  +2.52%  "
Set-TextValue "D33" "2.45"
Set-TextValue "E33" "  +5.54%  "
Set-TextValue "D34" "5.79"
Set-TextValue "E34" "  +5.57%  "
Set-TextValue "D35" "6.24"
Set-TextValue "E35" "  +7.83%  "
Set-TextValue "D36" "54.95"
Set-TextValue "E36" "  +0.64%  "
Set-TextValue "D37" "0.0412"
Set-TextValue "E37" "  +6.06%  "
Set-TextValue "D38" "448.68"
Set-TextValue "E38" "  +0.00%  "
Set-TextValue "D39" "0.0815"
Set-TextValue "E39" "  +1.73%  "
Set-TextValue "D40" "2.85"
Set-TextValue "E40" "  +11.80%  "
Set-TextValue "D41" "3.021.70"
Set-TextValue "E41" "  +3.39%  "
Set-TextValue "D42" "8.28"
Set-TextValue "D43" "0.116"
Set-TextValue "E43" "  +2.16%  "
Set-TextValue "D44" "27.89"
Set-TextValue "E44" "  +5.16%  "
Set-TextValue "D45" "0.265"
Set-TextValue "E45" "  +8.58%  "
Set-TextValue "D46" "2.22"
Set-TextValue "E46" "  +12.35%  "
Set-TextValue "E47" "  -0.02%  "
Set-TextValue "E48" "  +4.42%  "
Set-TextValue "D49" "118.21"
Set-TextValue "E49" "  +3.43%  "
Set-TextValue "D50" "0.0₃0517"
Set-TextValue "E50" "  +4.49%  "
Set-TextValue "D51" "2.08"
Set-TextValue "E51" "  +5.30%  "

Write-Output "Applied 86 cell updates"
